# Rename the default sheet from "Sheet1" to "Template".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Template"

# Populate the COA analyte/result table.
# Values are written in this order so the shared-string table comes out as
# Analyte(0), CBD(1), {{ cbd }}(2), Result(3).
$ws.Range("B5").Value = "Analyte"
$ws.Range("B8").Value = "CBD"
$ws.Range("C8").Value = "{{ cbd }}"
$ws.Range("C5").Value = "Result"

# Widen columns C and D to fit the new content.
$ws.Columns.Item(3).ColumnWidth = 22.5
$ws.Columns.Item(4).ColumnWidth = 19

# Leave the active selection on the result placeholder cell.
$ws.Range("C8").Select() | Out-Null
